$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.461.05'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.567.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.09'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0592'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.790.28'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.576.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.518'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.460.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '213.49'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('E24').Value = '  +2.79%  '
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.69'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.11%  '
$ws.Range('E29').Value = '  -1.85%  '
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.375.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.89%  '
$ws.Range('E34').Value = '  +1.17%  '
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.955'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.534'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.87%  '
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.702.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.43'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0957'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('E51').Value = '  -0.54%  '
